$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ebi3"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.928253333333333
$ws.Range("H2").Value = 17.78476
$ws.Range("I2").Value = 0.4729671618337555
$ws.Range("J2").Value = 0.4729671618337555
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.51056
$ws.Range("N2").Value = 4.53168
$ws.Range("O2").Value = 0.4569399709189402
$ws.Range("P2").Value = 0.4674407903484196
$ws.Range("Q2").Value = 8.954982355199999
$ws.Range("R2").Value = 80.59484119679999
$ws.Range("S2").Value = 0.2161176011739299
$ws.Range("T2").Value = 0.2210841439364195

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ebi3"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.928253333333333
$ws.Range("H3").Value = 17.78476
$ws.Range("I3").Value = 0.4729671618337555
$ws.Range("J3").Value = 0.4729671618337555
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9101323333333333
$ws.Range("N3").Value = 2.730397
$ws.Range("O3").Value = 0.2753123622535487
$ws.Range("P3").Value = 0.281639244528509
$ws.Range("Q3").Value = 5.395495038857778
$ws.Range("R3").Value = 48.55945534972
$ws.Range("S3").Value = 0.1302137065928077
$ws.Range("T3").Value = 0.133206114145652

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ebi3"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.928253333333333
$ws.Range("H4").Value = 17.78476
$ws.Range("I4").Value = 0.4729671618337555
$ws.Range("J4").Value = 0.4729671618337555
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.250344
$ws.Range("N4").Value = 0.751032
$ws.Range("O4").Value = 0.07572832597164705
$ws.Range("P4").Value = 0.07746861906775286
$ws.Range("Q4").Value = 1.48410265248
$ws.Range("R4").Value = 13.35692387232
$ws.Range("S4").Value = 0.03581701140523138
$ws.Range("T4").Value = 0.03664011289165543

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Ebi3"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.928253333333333
$ws.Range("H5").Value = 17.78476
$ws.Range("I5").Value = 0.4729671618337555
$ws.Range("J5").Value = 0.4729671618337555
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4119903333333333
$ws.Range("N5").Value = 1.235971
$ws.Range("O5").Value = 0.1246258678451818
$ws.Range("P5").Value = 0.1274898627192844
$ws.Range("Q5").Value = 2.442383066884444
$ws.Range("R5").Value = 21.98144760196
$ws.Range("S5").Value = 0.05894394300580431
$ws.Range("T5").Value = 0.06029851853291504

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ebi3"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.928253333333333
$ws.Range("H6").Value = 17.78476
$ws.Range("I6").Value = 0.4729671618337555
$ws.Range("J6").Value = 0.4729671618337555
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2227905
$ws.Range("N6").Value = 0.445581
$ws.Range("O6").Value = 0.06739347301068223
$ws.Range("P6").Value = 0.04596148333603414
$ws.Range("Q6").Value = 1.32075852426
$ws.Range("R6").Value = 7.92455114556
$ws.Range("S6").Value = 0.03187489965598218
$ws.Range("T6").Value = 0.02173827232711352

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ebi3"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.605922
$ws.Range("H7").Value = 19.817766
$ws.Range("I7").Value = 0.5270328381662445
$ws.Range("J7").Value = 0.5270328381662445
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.51056
$ws.Range("N7").Value = 4.53168
$ws.Range("O7").Value = 0.4569399709189402
$ws.Range("P7").Value = 0.4674407903484196
$ws.Range("Q7").Value = 9.97864153632
$ws.Range("R7").Value = 89.80777382688
$ws.Range("S7").Value = 0.2408223697450103
$ws.Range("T7").Value = 0.246356646412

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ebi3"
$ws.Range("C8").Value = "Il27ra"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.605922
$ws.Range("H8").Value = 19.817766
$ws.Range("I8").Value = 0.5270328381662445
$ws.Range("J8").Value = 0.5270328381662445
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9101323333333333
$ws.Range("N8").Value = 2.730397
$ws.Range("O8").Value = 0.2753123622535487
$ws.Range("P8").Value = 0.281639244528509
$ws.Range("Q8").Value = 6.012263203678
$ws.Range("R8").Value = 54.11036883310199
$ws.Range("S8").Value = 0.145098655660741
$ws.Range("T8").Value = 0.148433130382857

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ebi3"
$ws.Range("C9").Value = "Il27ra"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.605922
$ws.Range("H9").Value = 19.817766
$ws.Range("I9").Value = 0.5270328381662445
$ws.Range("J9").Value = 0.5270328381662445
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.250344
$ws.Range("N9").Value = 0.751032
$ws.Range("O9").Value = 0.07572832597164705
$ws.Range("P9").Value = 0.07746861906775286
$ws.Range("Q9").Value = 1.653752937168
$ws.Range("R9").Value = 14.883776434512
$ws.Range("S9").Value = 0.03991131456641567
$ws.Range("T9").Value = 0.04082850617609744

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Ebi3"
$ws.Range("C10").Value = "Il27ra"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.605922
$ws.Range("H10").Value = 19.817766
$ws.Range("I10").Value = 0.5270328381662445
$ws.Range("J10").Value = 0.5270328381662445
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4119903333333333
$ws.Range("N10").Value = 1.235971
$ws.Range("O10").Value = 0.1246258678451818
$ws.Range("P10").Value = 0.1274898627192844
$ws.Range("Q10").Value = 2.721576006753999
$ws.Range("R10").Value = 24.494184060786
$ws.Range("S10").Value = 0.06568192483937746
$ws.Range("T10").Value = 0.06719134418636932

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Ebi3"
$ws.Range("C11").Value = "Il27ra"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.605922
$ws.Range("H11").Value = 19.817766
$ws.Range("I11").Value = 0.5270328381662445
$ws.Range("J11").Value = 0.5270328381662445
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2227905
$ws.Range("N11").Value = 0.445581
$ws.Range("O11").Value = 0.06739347301068223
$ws.Range("P11").Value = 0.04596148333603414
$ws.Range("Q11").Value = 1.471736665341
$ws.Range("R11").Value = 8.830419992046
$ws.Range("S11").Value = 0.03551857335470005
$ws.Range("T11").Value = 0.02422321100892062
